$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "80×29=2320" "18×92=1656"
Replace-Text "25×22=550" "70×66=4620"
Replace-Text "89×97=8633" "70×67=4690"
Replace-Text "41×43=1763" "15×13=195"
Replace-Text "88×94=8272" "34×65=2210"
Replace-Text "12×49=588" "50×75=3750"
Replace-Text "45×66=2970" "62×57=3534"
Replace-Text "43×57=2451" "45×71=3195"
Replace-Text "39×52=2028" "55×64=3520"
Replace-Text "98×71=6958" "52×76=3952"
Replace-Text "40×35=1400" "70×80=5600"
Replace-Text "69×68=4692" "82×59=4838"
Replace-Text "65×98=6370" "90×94=8460"
Replace-Text "12×42=504" "67×37=2479"
Replace-Text "11×30=330" "48×40=1920"
Replace-Text "40×50=2000" "24×44=1056"
Replace-Text "20×47=940" "34×67=2278"
Replace-Text "50×42=2100" "92×14=1288"
Replace-Text "65×84=5460" "75×13=975"
Replace-Text "77×14=1078" "44×19=836"
Replace-Text "24×47=1128" "13×64=832"
Replace-Text "75×63=4725" "59×46=2714"
Replace-Text "56×49=2744" "73×48=3504"
Replace-Text "56×17=952" "75×39=2925"
Replace-Text "36×30=1080" "90×20=1800"
